$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.768.79'
$ws.Range("E2").Value = '  +0.75%  '
$ws.Range("D3").Value = '2.480.71'
$ws.Range("E3").Value = '  -0.12%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''587.31'
$ws.Range("E5").Value = '  +0.24%  '
$ws.Range("D6").Value = '''175.15'
$ws.Range("E6").Value = '  +1.80%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").Value = '''0.514'
$ws.Range("E8").Value = '  -0.50%  '
$ws.Range("D9").Value = '''0.145'
$ws.Range("E9").Value = '  +4.59%  '
$ws.Range("D10").Value = '''0.163'
$ws.Range("E10").Value = '  -1.53%  '
$ws.Range("D11").Value = '''4.96'
$ws.Range("E11").Value = '  +0.33%  '
$ws.Range("D12").Value = '''0.334'
$ws.Range("E12").Value = '  +0.17%  '
$ws.Range("D13").Value = '2.932.22'
$ws.Range("E13").Value = '  -0.15%  '
$ws.Range("D14").Value = '''25.28'
$ws.Range("E14").Value = '  -1.18%  '
$ws.Range("D15").Value = '67.723.45'
$ws.Range("E15").Value = '  +0.98%  '
$ws.Range("D16").Value = '''0.0000170'
$ws.Range("E16").Value = '  -0.34%  '
$ws.Range("D17").Value = '2.480.58'
$ws.Range("E17").Value = '  -0.17%  '
$ws.Range("B18").Value = 'Uniswap'
$ws.Range("C18").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D18").Value = '''7.44'
$ws.Range("E18").Value = '  -4.08%  '
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").Value = '''10.81'
$ws.Range("E19").Value = '  -2.06%  '
$ws.Range("D20").Value = '''348.62'
$ws.Range("E20").Value = '  -1.14%  '
$ws.Range("D21").Value = '''4.05'
$ws.Range("E21").Value = '  +1.28%  '
$ws.Range("D22").Value = '''1.00'
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").Value = '''70.66'
$ws.Range("E23").Value = '  +2.38%  '
$ws.Range("D24").Value = '''4.20'
$ws.Range("E24").Value = '  -1.19%  '
$ws.Range("D25").Value = '''1.69'
$ws.Range("E25").Value = '  -6.51%  '
$ws.Range("D26").Value = '''8.80'
$ws.Range("E26").Value = '  -5.35%  '
$ws.Range("D27").Value = '2.603.61'
$ws.Range("E27").Value = '  +1.12%  '
$ws.Range("D28").Value = '''0.997'
$ws.Range("E28").Value = '  -0.31%  '
$ws.Range("D29").Value = '0.0₃0892'
$ws.Range("E29").Value = '  -2.60%  '
$ws.Range("D30").Value = '''7.74'
$ws.Range("E30").Value = '  -0.06%  '
$ws.Range("D31").Value = '''495.95'
$ws.Range("E31").Value = '  -2.97%  '
$ws.Range("D32").Value = '''1.25'
$ws.Range("E32").Value = '  -0.10%  '
$ws.Range("D33").Value = '''1.76'
$ws.Range("E33").Value = '  -0.75%  '
$ws.Range("D35").Value = '''164.41'
$ws.Range("E35").Value = '  +1.45%  '
$ws.Range("D36").Value = '''0.120'
$ws.Range("E36").Value = '  +1.25%  '
$ws.Range("D38").Value = '''18.33'
$ws.Range("E38").Value = '  +0.78%  '
$ws.Range("E39").Value = '  -0.01%  '
$ws.Range("D40").Value = '''1.29'
$ws.Range("E40").Value = '  -3.39%  '
$ws.Range("D41").Value = '''1.73'
$ws.Range("E41").Value = '  +1.32%  '
$ws.Range("D42").Value = '''0.326'
$ws.Range("E42").Value = '  -1.27%  '
$ws.Range("D43").Value = '''4.78'
$ws.Range("E43").Value = '  -1.47%  '
$ws.Range("D44").Value = '''2.38'
$ws.Range("E44").Value = '  -0.44%  '
$ws.Range("D45").Value = '''147.79'
$ws.Range("E45").Value = '  +2.62%  '
$ws.Range("D46").Value = '''3.53'
$ws.Range("E46").Value = '  +0.70%  '
$ws.Range("D47").Value = '''0.510'
$ws.Range("E47").Value = '  -1.19%  '
$ws.Range("D48").Value = '0.0₆0254'
$ws.Range("E48").Value = '  -3.90%  '
$ws.Range("E49").Value = '  -0.10%  '
$ws.Range("D50").Value = '''1.56'
$ws.Range("E50").Value = '  -1.39%  '
$ws.Range("D51").Value = '''0.576'
$ws.Range("E51").Value = '  -1.63%  '
